$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire 4th row (previously ECs/Resolving-Mac row) which removes the need for
# the "Resolving-Mac" shared string as well since it was only used there.
$ws.Rows.Item(4).Delete()

# Update row 2 values (columns E through T) with the newly computed TPM-based figures.
$ws.Range("I2").Value = 0.7099439172299504
$ws.Range("J2").Value = 0.7099439172299504
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6124329999999999
$ws.Range("N2").Value = 1.837299
$ws.Range("Q2").Value = 0.2606200465726666
$ws.Range("R2").Value = 2.345580419154
$ws.Range("S2").Value = 0.7099439172299504
$ws.Range("T2").Value = 0.7099439172299504

# Update row 3 values (columns E through T) with the newly computed TPM-based figures.
$ws.Range("I3").Value = 0.2900560827700495
$ws.Range("J3").Value = 0.2900560827700495
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6124329999999999
$ws.Range("N3").Value = 1.837299
$ws.Range("Q3").Value = 0.106479438679
$ws.Range("R3").Value = 0.9583149481109998
$ws.Range("S3").Value = 0.2900560827700495
$ws.Range("T3").Value = 0.2900560827700495
